$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing tube IDs to the new "TESTImport" naming scheme and
# append five more rows so the list now runs from 0001 to 0010.
$ws.Range("A2").Value  = "TESTImport0001"
$ws.Range("A3").Value  = "TESTImport0002"
$ws.Range("A4").Value  = "TESTImport0003"
$ws.Range("A5").Value  = "TESTImport0004"
$ws.Range("A6").Value  = "TESTImport0005"
$ws.Range("A7").Value  = "TESTImport0006"
$ws.Range("A8").Value  = "TESTImport0007"
$ws.Range("A9").Value  = "TESTImport0008"
$ws.Range("A10").Value = "TESTImport0009"
$ws.Range("A11").Value = "TESTImport0010"

$ws.Range("A12").Select()
